$wb = $excel.ActiveWorkbook

# --- Update conversion text on "Hoja1" sheet ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$wsHoja1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 8.2 = 33853.28 pesos`n✅ 33853.28 pesos = 8.17 = 953.17 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

# --- Update rate figures on "tasas" sheet ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 122
$wsTasas.Range("O10").Value = 4130.1
$wsTasas.Range("N12").Value = 4143
$wsTasas.Range("O12").Value = 116.65
